$d = $word.ActiveDocument

$replacements = @(
    @("674÷3=224, 2", "704÷8=88, 0"),
    @("218÷9=24, 2", "523÷7=74, 5"),
    @("188÷7=26, 6", "935÷4=233, 3"),
    @("426÷4=106, 2", "391÷7=55, 6"),
    @("624÷7=89, 1", "439÷3=146, 1"),
    @("312÷9=34, 6", "943÷3=314, 1"),
    @("116÷4=29, 0", "743÷8=92, 7"),
    @("438÷3=146, 0", "533÷9=59, 2"),
    @("595÷3=198, 1", "200÷8=25, 0"),
    @("644÷5=128, 4", "607÷5=121, 2"),
    @("720÷4=180, 0", "234÷6=39, 0"),
    @("569÷6=94, 5", "134÷9=14, 8"),
    @("313÷5=62, 3", "476÷3=158, 2"),
    @("986÷5=197, 1", "178÷3=59, 1"),
    @("714÷9=79, 3", "579÷7=82, 5"),
    @("874÷2=437, 0", "723÷9=80, 3"),
    @("945÷5=189, 0", "925÷4=231, 1"),
    @("637÷9=70, 7", "280÷7=40, 0"),
    @("793÷3=264, 1", "165÷3=55, 0"),
    @("156÷6=26, 0", "117÷3=39, 0"),
    @("373÷8=46, 5", "196÷8=24, 4"),
    @("110÷2=55, 0", "886÷3=295, 1"),
    @("159÷4=39, 3", "235÷5=47, 0"),
    @("764÷7=109, 1", "911÷5=182, 1"),
    @("881÷5=176, 1", "485÷4=121, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
